# Add a new GMPE entry "NCREE2011" (fortran_name "NCREE_2011", region
# "Taiwan") to the GMPE list on the active worksheet. The new record is
# inserted as row 58, pushing the existing rows 58-69 down to 59-70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 58; Excel shifts rows 58:69 (and
# their formatting, e.g. the wrapped/tall row that used to be row 69)
# down to 59:70 automatically.
$ws.Rows.Item(58).Insert()

# Populate the new row with the NCREE2011 model data.
$ws.Range("A58").Value = "NCREE2011"
$ws.Range("B58").Value = "NCREE_2011"
$ws.Range("C58").Value = 1
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = "Taiwan"
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0.01
$ws.Range("L58").Value = 10

# Match the author's final view state: frozen header row, scrolled so
# row 2 is the first visible row below the freeze, with H58 selected.
[void]$ws.Range("H58").Select()
$excel.ActiveWindow.ScrollRow = 2
